$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ref_img" column (E) entirely - header + data no longer used
$ws.Range("E1:E4").ClearContents()

# Update task/video ids for existing rows 2 and 3
$ws.Range("D2").Value = "task_01k2sd81t6fnja114fgn43mydt"
$ws.Range("D3").Value = "task_01k2sda8fdfgg9hx4948a3p9st"

# Row 4 now becomes a different prompt/task (previously the "girl" reference-image prompt)
$ws.Range("A4").Value = "aerial shot of Victoria Waterfal, Zimbabwe"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "success"
$ws.Range("D4").Value = "task_01k2sdc5mpfk7rtgt397bq0ecy"

# New row 5
$ws.Range("A5").Value = "aerial video of Lumangwe waterfal"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "success"
$ws.Range("D5").Value = "task_01k2sddqq2e42rs2cefmypkryq"

# New row 6
$ws.Range("A6").Value = "aerial video of Singapore Marinabay"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "success"
$ws.Range("D6").Value = "task_01k2sdfqjbfe8szhjdqmba5xvq"
